$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage-looking text values need NumberFormat forced to text
# to prevent Excel auto-converting "NN%" strings into numeric percentages.
$percentCells = @("H4", "H6", "H7", "H9", "H12", "H15", "H17", "H18", "H19", "H20", "H21", "H22", "H23", "H25", "H27", "H28", "H30", "H31", "H34", "H38", "H39", "H40", "H46")
foreach ($cellRef in $percentCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("E2").Value = "2026-02-28 16:48:30"
$ws.Range("K2").Value = "11.6 MJ/m2"
$ws.Range("E3").Value = "2026-02-28 16:48:32"
$ws.Range("K3").Value = "11.2 MJ/m2"
$ws.Range("E4").Value = "2026-02-28 16:48:35"
$ws.Range("H4").Value = "82%"
$ws.Range("K4").Value = "5.9 MJ/m2"
$ws.Range("O4").Value = "11.0 °C"
$ws.Range("E5").Value = "2026-02-28 16:48:37"
$ws.Range("K5").Value = "8.6 MJ/m2"
$ws.Range("N5").Value = "-2.2 °C 16:29 TU"
$ws.Range("E6").Value = "2026-02-28 16:48:39"
$ws.Range("H6").Value = "84%"
$ws.Range("K6").Value = "10.1 MJ/m2"
$ws.Range("O6").Value = "12.1 °C"
$ws.Range("E7").Value = "2026-02-28 16:48:41"
$ws.Range("H7").Value = "76%"
$ws.Range("K7").Value = "5.7 MJ/m2"
$ws.Range("E8").Value = "2026-02-28 16:48:43"
$ws.Range("E9").Value = "2026-02-28 16:48:46"
$ws.Range("H9").Value = "83%"
$ws.Range("K9").Value = "11.1 MJ/m2"
$ws.Range("O9").Value = "11.2 °C"
$ws.Range("E10").Value = "2026-02-28 16:48:48"
$ws.Range("K10").Value = "8.5 MJ/m2"
$ws.Range("O10").Value = "11.0 °C"
$ws.Range("E11").Value = "2026-02-28 16:48:49"
$ws.Range("O11").Value = "6.8 °C"
$ws.Range("E12").Value = "2026-02-28 16:48:50"
$ws.Range("H12").Value = "85%"
$ws.Range("O12").Value = "10.6 °C"
$ws.Range("E13").Value = "2026-02-28 16:48:51"
$ws.Range("K13").Value = "11.9 MJ/m2"
$ws.Range("O13").Value = "6.3 °C"
$ws.Range("E14").Value = "2026-02-28 16:48:53"
$ws.Range("K14").Value = "5.6 MJ/m2"
$ws.Range("O14").Value = "12.7 °C"
$ws.Range("E15").Value = "2026-02-28 16:48:54"
$ws.Range("H15").Value = "81%"
$ws.Range("O15").Value = "11.1 °C"
$ws.Range("E16").Value = "2026-02-28 16:48:55"
$ws.Range("K16").Value = "11.3 MJ/m2"
$ws.Range("E17").Value = "2026-02-28 16:48:56"
$ws.Range("H17").Value = "78%"
$ws.Range("K17").Value = "10.5 MJ/m2"
$ws.Range("O17").Value = "3.1 °C"
$ws.Range("E18").Value = "2026-02-28 16:48:57"
$ws.Range("H18").Value = "83%"
$ws.Range("K18").Value = "9.6 MJ/m2"
$ws.Range("L18").Value = "20.5 km/h - 87º 16:22 TU"
$ws.Range("O18").Value = "11.7 °C"
$ws.Range("E19").Value = "2026-02-28 16:48:58"
$ws.Range("H19").Value = "75%"
$ws.Range("K19").Value = "7.6 MJ/m2"
$ws.Range("E20").Value = "2026-02-28 16:48:59"
$ws.Range("H20").Value = "55%"
$ws.Range("K20").Value = "13.8 MJ/m2"
$ws.Range("E21").Value = "2026-02-28 16:49:00"
$ws.Range("H21").Value = "70%"
$ws.Range("K21").Value = "10.8 MJ/m2"
$ws.Range("E22").Value = "2026-02-28 16:49:02"
$ws.Range("H22").Value = "62%"
$ws.Range("K22").Value = "11.3 MJ/m2"
$ws.Range("E23").Value = "2026-02-28 16:49:05"
$ws.Range("H23").Value = "66%"
$ws.Range("K23").Value = "12.1 MJ/m2"
$ws.Range("E24").Value = "2026-02-28 16:49:08"
$ws.Range("K24").Value = "2.6 MJ/m2"
$ws.Range("O24").Value = "8.2 °C"
$ws.Range("E25").Value = "2026-02-28 16:49:11"
$ws.Range("H25").Value = "55%"
$ws.Range("K25").Value = "11.6 MJ/m2"
$ws.Range("E26").Value = "2026-02-28 16:49:13"
$ws.Range("O26").Value = "5.1 °C"
$ws.Range("E27").Value = "2026-02-28 16:49:16"
$ws.Range("H27").Value = "47%"
$ws.Range("K27").Value = "13.7 MJ/m2"
$ws.Range("E28").Value = "2026-02-28 16:49:19"
$ws.Range("H28").Value = "83%"
$ws.Range("K28").Value = "7.0 MJ/m2"
$ws.Range("O28").Value = "9.3 °C"
$ws.Range("E29").Value = "2026-02-28 16:49:21"
$ws.Range("K29").Value = "12.2 MJ/m2"
$ws.Range("O29").Value = "11.7 °C"
$ws.Range("E30").Value = "2026-02-28 16:49:24"
$ws.Range("H30").Value = "81%"
$ws.Range("J30").Value = "1024.7 hPa"
$ws.Range("K30").Value = "12.5 MJ/m2"
$ws.Range("O30").Value = "10.9 °C"
$ws.Range("E31").Value = "2026-02-28 16:49:26"
$ws.Range("H31").Value = "84%"
$ws.Range("K31").Value = "11.8 MJ/m2"
$ws.Range("L31").Value = "53.6 km/h - 337º 16:29 TU"
$ws.Range("O31").Value = "11.6 °C"
$ws.Range("E32").Value = "2026-02-28 16:49:29"
$ws.Range("K32").Value = "2.8 MJ/m2"
$ws.Range("E33").Value = "2026-02-28 16:49:32"
$ws.Range("K33").Value = "11.6 MJ/m2"
$ws.Range("E34").Value = "2026-02-28 16:49:34"
$ws.Range("H34").Value = "61%"
$ws.Range("K34").Value = "11.0 MJ/m2"
$ws.Range("E35").Value = "2026-02-28 16:49:36"
$ws.Range("K35").Value = "5.2 MJ/m2"
$ws.Range("E36").Value = "2026-02-28 16:49:38"
$ws.Range("J36").Value = "1024.9 hPa"
$ws.Range("K36").Value = "13.4 MJ/m2"
$ws.Range("O36").Value = "12.4 °C"
$ws.Range("E37").Value = "2026-02-28 16:49:41"
$ws.Range("O37").Value = "6.8 °C"
$ws.Range("E38").Value = "2026-02-28 16:49:44"
$ws.Range("H38").Value = "80%"
$ws.Range("K38").Value = "7.7 MJ/m2"
$ws.Range("O38").Value = "11.7 °C"
$ws.Range("E39").Value = "2026-02-28 16:49:46"
$ws.Range("H39").Value = "56%"
$ws.Range("K39").Value = "14.6 MJ/m2"
$ws.Range("O39").Value = "-0.3 °C"
$ws.Range("E40").Value = "2026-02-28 16:49:49"
$ws.Range("H40").Value = "77%"
$ws.Range("J40").Value = "1024.5 hPa"
$ws.Range("O40").Value = "7.1 °C"
$ws.Range("E41").Value = "2026-02-28 16:49:51"
$ws.Range("K41").Value = "5.9 MJ/m2"
$ws.Range("E42").Value = "2026-02-28 16:49:54"
$ws.Range("O42").Value = "11.0 °C"
$ws.Range("E43").Value = "2026-02-28 16:49:57"
$ws.Range("K43").Value = "6.3 MJ/m2"
$ws.Range("O43").Value = "7.0 °C"
$ws.Range("E44").Value = "2026-02-28 16:49:59"
$ws.Range("K44").Value = "12.9 MJ/m2"
$ws.Range("E45").Value = "2026-02-28 16:50:02"
$ws.Range("E46").Value = "2026-02-28 16:50:04"
$ws.Range("H46").Value = "77%"
$ws.Range("K46").Value = "4.6 MJ/m2"
$ws.Range("O46").Value = "11.6 °C"
